$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 6 new WATIAM test-case rows (rows 3-8), matching the formatting of the
# existing template row (row 2).
for ($r = 3; $r -le 8; $r++) {
    $ws.Range("A2:E2").Copy()
    $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# Fill in the new values in the same order they were originally authored.
$ws.Range("A3").Value = "WATIAM002"
$ws.Range("B3").Value = "WAT-226"
$ws.Range("C3").Value = "Verify that the app name is displayed as SaR Labs"
$ws.Range("D3").Value = "Y"

$ws.Range("B4").Value = "WAT-227"
$ws.Range("C4").Value = 'Verify that the contact email is display as " sarlabs.info@clarivate.com " '
$ws.Range("A4").Value = "WATIAM003"
$ws.Range("D4").Value = "Y"

$ws.Range("A5").Value = "WATIAM004"
$ws.Range("B5").Value = "WAT-228"
$ws.Range("C5").Value = "Verify that the marketing text module should not appear in the login page"
$ws.Range("D5").Value = "Y"

$ws.Range("A6").Value = "WATIAM005"
$ws.Range("B6").Value = "WAT-229"
$ws.Range("C6").Value = "Verify that the links to the Terms and Privacy Statement should appear in login page"
$ws.Range("D6").Value = "Y"

$ws.Range("A7").Value = "WATIAM006"
$ws.Range("C7").Value = "Verify that Link to forgot password should be available on the sign in page"
$ws.Range("B7").Value = "WAT-232"
$ws.Range("D7").Value = "Y"

$ws.Range("A8").Value = "WATIAM007"
$ws.Range("C8").Value = "Verify that standard platform auth error messaging is displayed when Email address is not properly formatted"
$ws.Range("B8").Value = "WAT-235"
$ws.Range("D8").Value = "Y"

$ws.Range("B8").Select()
